$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AdminUserName (H2) and Id (K2) values
$ws.Range("H2").Value = "rajakolla3409"
$ws.Range("K2").Value = "'0247"

# Update selected/active cell to H2
$ws.Range("H2").Select()
